# Denmark Division 2 - base update (28-05-2024 07:50)
# The underlying data source re-sorted a handful of matches that share the
# same matchday/date. For each group below, the full record (every column
# except the running index in column A) moves to a different row while the
# row's own index in column A stays put.
#
# Groups (rows that exchange their B:AD data):
#   30 <-> 31                (swap)
#   70 -> 71 -> 72 -> 70      (rotate)
#   130 <-> 131               (swap)
#   157 -> 158 -> 159 -> 157  (rotate)
#   174 <-> 176                (swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 2")

$colFrom = "B"
$colTo   = "AD"

function Get-RowData($rowNum) {
    $rng = $ws.Range("$colFrom$rowNum`:$colTo$rowNum")
    return $rng.Value2
}

function Set-RowData($rowNum, $data) {
    $rng = $ws.Range("$colFrom$rowNum`:$colTo$rowNum")
    $rng.Value2 = $data
}

# Capture the current ("before") data for every row that participates in a shuffle.
$rows = 30,31,70,71,72,130,131,157,158,159,174,176
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = Get-RowData $r
}

# Mapping: destination row -> source row (whose B:AD values land there)
$mapping = @{
    30  = 31
    31  = 30
    70  = 72
    71  = 70
    72  = 71
    130 = 131
    131 = 130
    157 = 158
    158 = 159
    159 = 157
    174 = 176
    176 = 174
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    Set-RowData $dest $orig[$src]
}
